$d = $word.ActiveDocument

$pairs = @(
    @("31×48=", "75×19="),
    @("96×79=", "63×33="),
    @("46×81=", "42×74="),
    @("93×97=", "48×56="),
    @("48×28=", "22×14="),
    @("45×45=", "23×81="),
    @("14×24=", "47×45="),
    @("90×16=", "61×16="),
    @("15×70=", "95×45="),
    @("23×61=", "29×92="),
    @("27×16=", "23×14="),
    @("48×45=", "73×90="),
    @("70×73=", "91×85="),
    @("59×22=", "72×14="),
    @("76×73=", "86×83="),
    @("62×72=", "80×72="),
    @("17×11=", "69×91="),
    @("91×59=", "98×98="),
    @("56×44=", "61×38="),
    @("45×62=", "36×60="),
    @("66×14=", "43×99="),
    @("82×38=", "73×84="),
    @("93×96=", "65×37="),
    @("73×99=", "14×61="),
    @("59×61=", "56×12=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
